$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "register register" opcode table (rows 26-35) had a duplicated entry:
# row 35 ("str ", F6, "00 r r r r r r") duplicated the opcode sequencing used
# by the following "register register indirect 16bits" table. Remove row 35's
# content entirely (first 4 addressing modes done; this one still needs work).
$ws.Range("A35:C35").ClearContents()

# Column B (opcode hex) for the remaining two tables shifts up by one slot,
# since the old B35 value ("F6") now belongs at the top of the
# "register register indirect 16bits" table (B38), and the table's last
# value ("E3") drops off the end.
$ws.Range("B38").Value = "F6"
$ws.Range("B39").Value = "F5"
$ws.Range("B40").Value = "F4"
$ws.Range("B41").Value = "F3"
$ws.Range("B42").Value = "F2"
$ws.Range("B43").Value = "F1"
$ws.Range("B44").Value = "F0"
$ws.Range("B45").Value = "EF"
$ws.Range("B46").Value = "EE"
$ws.Range("B47").Value = "ED"

$ws.Range("B50").Value = "EC"
$ws.Range("B51").Value = "EB"
$ws.Range("B52").Value = "EA"
$ws.Range("B53").Value = "E9"
$ws.Range("B54").Value = "E8"
$ws.Range("B55").Value = "E7"
$ws.Range("B56").Value = "E6"
$ws.Range("B57").Value = "E5"
$ws.Range("B58").Value = "E4"

# Restore the view to where the user left off editing.
$ws.Range("C35").Select()
$ws.Application.ActiveWindow.ScrollRow = 22
